$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.288.51'
$ws.Range("E2").Value = '  +2.08%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.944.62'
$ws.Range("E3").Value = '  +2.12%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.38'
$ws.Range("E5").Value = '  +1.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.87'
$ws.Range("E6").Value = '  +5.56%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.945.33'
$ws.Range("E7").Value = '  +2.17%  '

# Row 8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.540'
$ws.Range("E9").Value = '  +1.79%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.172'
$ws.Range("E10").Value = '  +2.73%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.47'
$ws.Range("E11").Value = '  +2.49%  '

# Row 12
$ws.Range("E12").Value = '  +3.13%  '

# Row 13
$ws.Range("E13").Value = '  +6.36%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.70'
$ws.Range("E14").Value = '  +5.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.612.83'
$ws.Range("E15").Value = '  +2.27%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.933.55'
$ws.Range("E16").Value = '  +1.74%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.114.07'
$ws.Range("E17").Value = '  +1.57%  '

# Row 18
$ws.Range("E18").Value = '  +2.31%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.41'
$ws.Range("E19").Value = '  +7.53%  '

# Row 20
$ws.Range("E20").Value = '  -0.85%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.10'
$ws.Range("E21").Value = '  -2.25%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '500.11'
$ws.Range("E22").Value = '  +3.06%  '

# Row 23
$ws.Range("E23").Value = '  +4.03%  '

# Row 24
$ws.Range("E24").Value = '  +5.90%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.20'
$ws.Range("E25").Value = '  +2.65%  '

# Row 26
$ws.Range("E26").Value = '  +3.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.44'
$ws.Range("E27").Value = '  +2.60%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.31'
$ws.Range("E28").Value = '  +2.69%  '

# Row 29
$ws.Range("E29").Value = '  +0.16%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.03'
$ws.Range("E30").Value = '  +1.81%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.46'
$ws.Range("E31").Value = '  +3.04%  '

# Row 32
$ws.Range("B32").Value = 'WrappedeETH'
$ws.Range("C32").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.099.40'
$ws.Range("E32").Value = '  +2.15%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.92'
$ws.Range("E33").Value = '  -0.54%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '32.45'
$ws.Range("E34").Value = '  +0.50%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.910.95'
$ws.Range("E35").Value = '  +2.54%  '

# Row 36
$ws.Range("E36").Value = '  +2.30%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.21'
$ws.Range("E37").Value = '  +5.09%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.04'
$ws.Range("E38").Value = '  +0.70%  '

# Row 39
$ws.Range("E39").Value = '  +1.51%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.29'
$ws.Range("E40").Value = '  +10.10%  '

# Row 41
$ws.Range("E41").Value = '  +3.74%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.05%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.14'
$ws.Range("E43").Value = '  +8.01%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '442.90'
$ws.Range("E44").Value = '  +0.02%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.32'
$ws.Range("E45").Value = '  -0.34%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.72'
$ws.Range("E46").Value = '  +3.84%  '

# Row 47
$ws.Range("E47").Value = '  +0.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000282'
$ws.Range("E48").Value = '  +25.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0371'
$ws.Range("E49").Value = '  +3.82%  '

# Row 50
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '41.06'
$ws.Range("E50").Value = '  +6.55%  '

# Row 51
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.78'
$ws.Range("E51").Value = '  +0.35%  '
